$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 8 new arrival rows (rows 151-158) for "Saturday, Jan 14", mirroring the
# existing table layout: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=FROM, F=SHORT,
# G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank), L=DIFFERENCE, M=(blank)

# Row 151
$ws.Cells.Item(151, 1).Value = 150
$ws.Cells.Item(151, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(151, 3).Value = "6:19 PM"
$ws.Cells.Item(151, 4).Value = "UNKNOWN"
$ws.Cells.Item(151, 5).Value = "Funchal"
$ws.Cells.Item(151, 6).Value = "(FNC)"
$ws.Cells.Item(151, 7).Value = "Enter Air "
$ws.Cells.Item(151, 8).Value = "B738"
$ws.Cells.Item(151, 9).Value = "(SP-ESF)"
$ws.Cells.Item(151, 10).Value = "6:19 PM"
$ws.Cells.Item(151, 12).Value = "0 hours, 0 minutes"

# Row 152
$ws.Cells.Item(152, 1).Value = 151
$ws.Cells.Item(152, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(152, 3).Value = "6:45 PM"
$ws.Cells.Item(152, 4).Value = "FR3727"
$ws.Cells.Item(152, 5).Value = "Billund"
$ws.Cells.Item(152, 6).Value = "(BLL)"
$ws.Cells.Item(152, 7).Value = "Ryanair "
$ws.Cells.Item(152, 8).Value = "B738"
$ws.Cells.Item(152, 9).Value = "(SP-RSM)"
$ws.Cells.Item(152, 10).Value = "6:54 PM"
$ws.Cells.Item(152, 12).Value = "0 hours, 9 minutes"

# Row 153
$ws.Cells.Item(153, 1).Value = 152
$ws.Cells.Item(153, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(153, 3).Value = "9:30 PM"
$ws.Cells.Item(153, 4).Value = "FR1750"
$ws.Cells.Item(153, 5).Value = "London"
$ws.Cells.Item(153, 6).Value = "(STN)"
$ws.Cells.Item(153, 7).Value = "Ryanair "
$ws.Cells.Item(153, 8).Value = "B38M"
$ws.Cells.Item(153, 9).Value = "(EI-HEV)"
$ws.Cells.Item(153, 10).Value = "9:27 PM"
$ws.Cells.Item(153, 12).Value = "0 hours, -3 minutes"

# Row 154
$ws.Cells.Item(154, 1).Value = 153
$ws.Cells.Item(154, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(154, 3).Value = "10:05 PM"
$ws.Cells.Item(154, 4).Value = "FR9628"
$ws.Cells.Item(154, 5).Value = "Liverpool"
$ws.Cells.Item(154, 6).Value = "(LPL)"
$ws.Cells.Item(154, 7).Value = "Ryanair "
$ws.Cells.Item(154, 8).Value = "B738"
$ws.Cells.Item(154, 9).Value = "(EI-EKZ)"
$ws.Cells.Item(154, 10).Value = "9:56 PM"
$ws.Cells.Item(154, 12).Value = "0 hours, -9 minutes"

# Row 155
$ws.Cells.Item(155, 1).Value = 154
$ws.Cells.Item(155, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(155, 3).Value = "10:10 PM"
$ws.Cells.Item(155, 4).Value = "FR4271"
$ws.Cells.Item(155, 5).Value = "Budapest"
$ws.Cells.Item(155, 6).Value = "(BUD)"
$ws.Cells.Item(155, 7).Value = "Ryanair "
$ws.Cells.Item(155, 8).Value = "B738"
$ws.Cells.Item(155, 9).Value = "(SP-RSM)"
$ws.Cells.Item(155, 10).Value = "10:21 PM"
$ws.Cells.Item(155, 12).Value = "0 hours, 11 minutes"

# Row 156
$ws.Cells.Item(156, 1).Value = 155
$ws.Cells.Item(156, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(156, 3).Value = "10:10 PM"
$ws.Cells.Item(156, 4).Value = "FR7947"
$ws.Cells.Item(156, 5).Value = "Oslo"
$ws.Cells.Item(156, 6).Value = "(TRF)"
$ws.Cells.Item(156, 7).Value = "Ryanair "
$ws.Cells.Item(156, 8).Value = "B738"
$ws.Cells.Item(156, 9).Value = "(SP-RSX)"
$ws.Cells.Item(156, 10).Value = "11:24 PM"
$ws.Cells.Item(156, 12).Value = "1 hours, 14 minutes"

# Row 157
$ws.Cells.Item(157, 1).Value = 156
$ws.Cells.Item(157, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(157, 3).Value = "10:15 PM"
$ws.Cells.Item(157, 4).Value = "LH1380"
$ws.Cells.Item(157, 5).Value = "Frankfurt"
$ws.Cells.Item(157, 6).Value = "(FRA)"
$ws.Cells.Item(157, 7).Value = "Lufthansa "
$ws.Cells.Item(157, 8).Value = "CRJ9"
$ws.Cells.Item(157, 9).Value = "(D-ACNB)"
$ws.Cells.Item(157, 10).Value = "10:07 PM"
$ws.Cells.Item(157, 12).Value = "0 hours, -8 minutes"

# Row 158
$ws.Cells.Item(158, 1).Value = 157
$ws.Cells.Item(158, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(158, 3).Value = "11:10 PM"
$ws.Cells.Item(158, 4).Value = "FR9627"
$ws.Cells.Item(158, 5).Value = "Tel Aviv"
$ws.Cells.Item(158, 6).Value = "(TLV)"
$ws.Cells.Item(158, 7).Value = "Ryanair "
$ws.Cells.Item(158, 8).Value = "B738"
$ws.Cells.Item(158, 9).Value = "(SP-RKR)"
$ws.Cells.Item(158, 10).Value = "11:19 PM"
$ws.Cells.Item(158, 12).Value = "0 hours, 9 minutes"

